# "Correccion faltante 09 - Vera, Bonader y Diaz Perdomo TERMINADO OK"
#
# Marks the two pending groups (row 15 = group 10: Bonader/Vera and
# row 21 = group 16: Diaz/Rodriguez) as finished, the same way the other
# already-graded rows (6-11, 14, 17-19) are marked: a green fill on the
# group-number/name/partner cells and an "OK" note in column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$green = 5296274   # RGB(146, 208, 80) == FF92D050, same green already used on rows 6-14,17-19

# --- Row 15 (Grupo 10: Araceli Bonader / Lorenzo Vera) ---
$ws.Range("B15:D15").Interior.Color = $green
$ws.Range("E15").Value = "OK"

# --- Row 21 (Grupo 16: Facundo Diaz / David Rodriguez) ---
$ws.Range("B21:C21").Interior.Color = $green
$ws.Range("E21").Value = "OK"

# Leave the cursor/selection where the author left it after the edit
[void]$ws.Range("D6").Select()
